$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Community "size" (count) column corrections
$ws.Range("B2").Value = 3
$ws.Range("B5").Value = 0

# Corrected member lists per community (column C), written in row order so any
# freshly-introduced text is appended to the shared-string pool in that order.
$ws.Range("C2").Value = "['Fc5', 'Fc3', 'C5', 'F7', 'F5', 'F3', 'Ft7', 'T7', 'T9', 'Tp7']"
$ws.Range("C3").Value = "['Fc1', 'Fp1', 'Fpz', 'Fp2', 'Af7', 'Af3', 'Af4', 'Af8', 'T10', 'O1', 'O2', 'Iz']"
$ws.Range("C4").Value = "['Fcz', 'Fc2', 'Cz', 'Afz', 'F1', 'Fz', 'F2', 'Poz']"
$ws.Range("C5").Value = "['Fc4', 'Fc6', 'C2', 'C4', 'C6', 'Cp2', 'Cp4', 'Cp6', 'F4', 'F6', 'F8', 'Ft8', 'T8', 'Tp8', 'P2', 'P4', 'P6', 'P8', 'Po4', 'Po8', 'Oz']"
$ws.Range("C6").Value = "['C3', 'C1', 'Cp5', 'Cp3', 'Cp1', 'Cpz', 'P7', 'P5', 'P3', 'P1', 'Pz', 'Po7', 'Po3']"
